$wb = $excel.ActiveWorkbook

# ----- Sheet "About": drop the old "Update for Canada" note block (rows 31-35) -----
$s1 = $wb.Worksheets.Item("About")
$s1.Rows("31:35").Delete()

# ----- Sheet "BPaFF-BITPTaP": re-work plant-type list -----
$s2 = $wb.Worksheets.Item("BPaFF-BITPTaP")
# row 13 used to be "coal to gas" -> now "lignite"
$s2.Range("A13").Value2 = "lignite"
# three new rows appended, mirroring the flagged value of a related plant type
$s2.Range("A15").Value2 = "crude oil"
$s2.Range("B15").Formula = "=B11"
$s2.Range("A16").Value2 = "heavy or residual fuel oil"
$s2.Range("B16").Formula = "=B11"
$s2.Range("A17").Value2 = "municipal solid waste"
$s2.Range("B17").Formula = "=B9"
# header cell "Boolean" is now right aligned
$s2.Range("B1").HorizontalAlignment = -4152

# ----- Sheet "BPaFF-BDTPTPF": same restructuring, plus hydro flag flips to 0 -----
$s3 = $wb.Worksheets.Item("BPaFF-BDTPTPF")
$s3.Range("B5").Value2 = 0
$s3.Range("A13").Value2 = "lignite"
$s3.Range("A15").Value2 = "crude oil"
$s3.Range("B15").Formula = "=B11"
$s3.Range("A16").Value2 = "heavy or residual fuel oil"
$s3.Range("B16").Formula = "=B11"
$s3.Range("A17").Value2 = "municipal solid waste"
$s3.Range("B17").Formula = "=B9"
$s3.Range("B1").HorizontalAlignment = -4152

Write-Host "edits applied"
